$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    19  = @("sd", "Statement-non-opinion")
    22  = @("b", "Acknowledge (Backchannel)")
    32  = @("ba", "Appreciation")
    55  = @("b", "Acknowledge (Backchannel)")
    64  = @("b", "Acknowledge (Backchannel)")
    73  = @("aa", "Agree/Accept")
    88  = @("sd", "Statement-non-opinion")
    90  = @("sd", "Statement-non-opinion")
    92  = @("sd", "Statement-non-opinion")
    93  = @("sd", "Statement-non-opinion")
    125 = @("aa", "Agree/Accept")
    129 = @("aa", "Agree/Accept")
    139 = @("aa", "Agree/Accept")
    142 = @("b", "Acknowledge (Backchannel)")
    165 = @("b", "Acknowledge (Backchannel)")
    171 = @("b", "Acknowledge (Backchannel)")
    186 = @("sd", "Statement-non-opinion")
    192 = @("sv", "Statement-opinion")
    193 = @("aa", "Agree/Accept")
    205 = @("sd", "Statement-non-opinion")
    206 = @("b", "Acknowledge (Backchannel)")
    217 = @("sd", "Statement-non-opinion")
    220 = @("aa", "Agree/Accept")
    229 = @("sd", "Statement-non-opinion")
    234 = @("sd", "Statement-non-opinion")
    244 = @("b", "Acknowledge (Backchannel)")
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
